$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "username"
$ws.Range("C1").Value = "user_id"
$ws.Range("D1").Value = "created_at"

# Row 2 data
$ws.Range("A2").Value = "Umrbek Xudayorovich"
$ws.Range("B2").Value = "pipcoder"
$ws.Range("C2").Value = 324304236
$ws.Range("D2").Value = "2024-11-08T09:40:52.003144Z"
